$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "D999"
$ws.Range("B9").Value = "6deba86e6bbcba88fd88348250d93153"
$ws.Range("C9").Value = "08560d220c5cfe9fe0e86abbb69fc069f3337e0adece119ffb55d63a01cb2e9a"
